$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '39.534.39'
$ws.Range('E2').Value = '  +2.07%  '
$ws.Range('D3').Value = '2.164.77'
$ws.Range('E3').Value = '  +3.08%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '228.52'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.39%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.634'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +2.97%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '63.59'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +2.10%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  +1.21%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0853'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +1.63%  '
$ws.Range('E11').Value = '  +0.38%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '16.06'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +2.26%  '
$ws.Range('D13').Value = '2.485.85'
$ws.Range('E13').Value = '  +3.02%  '
$ws.Range('E14').Value = '  +0.07%  '
$ws.Range('E15').Value = '  +0.27%  '
$ws.Range('E16').Value = '  -0.28%  '
$ws.Range('D17').Value = '2.164.44'
$ws.Range('E17').Value = '  +2.74%  '
$ws.Range('D18').Value = '39.490.41'
$ws.Range('E18').Value = '  +2.07%  '
$ws.Range('E19').Value = '  +1.60%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '72.23'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.75%  '
$ws.Range('D21').Value = '0.0₃0847'
$ws.Range('E21').Value = '  +0.96%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '229.95'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.90%  '
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('E24').Value = '  -1.04%  '
$ws.Range('E25').Value = '  +1.62%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.69'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.67%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '172.10'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.15%  '
$ws.Range('E28').Value = '  +0.75%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '19.89'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +2.97%  '
$ws.Range('E30').Value = '  +1.11%  '
$ws.Range('E31').Value = '  +5.68%  '
$ws.Range('E32').Value = '  +2.12%  '
$ws.Range('E33').Value = '  +1.44%  '
$ws.Range('E34').Value = '  -0.63%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '7.03'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.19%  '
$ws.Range('E36').Value = '  +0.10%  '
$ws.Range('B37').Value = 'LidoDAOToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.44'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +2.12%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.69'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +4.32%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.999'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -0.35%  '
$ws.Range('B40').Value = 'Aave'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '102.63'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -0.11%  '
$ws.Range('B41').Value = 'InjectiveProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '17.92'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -1.07%  '
$ws.Range('E42').Value = '  +0.03%  '
$ws.Range('B43').Value = 'FTXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '4.62'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +11.75%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '1.523.40'
$ws.Range('E44').Value = '  -0.34%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.21'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.32%  '
$ws.Range('E46').Value = '  +5.30%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0924'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +1.60%  '
$ws.Range('B48').Value = 'HuobiToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.81'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.07%  '
$ws.Range('E49').Value = '  -0.81%  '
$ws.Range('D50').Value = '2.369.24'
$ws.Range('E50').Value = '  +3.05%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.96'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.35%  '
